$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the desc text of row 32 (胡萝卜须) - content of existing shared string changes in-place
$ws.Range("I32").Value = "饱腹:{0}\n治疗:{1}\n最大血量:{2}\n"

# Rows 33-42: desc column (I) text updates (food buff descriptions)
$ws.Range("I33").Value = "饱腹:{0}\n治疗:{1}\n力量:{2}\n"
$ws.Range("I34").Value = "饱腹:{0}\n治疗:{1}\n力量:{2}\n"
$ws.Range("I35").Value = "饱腹:{0}\n治疗:{1}\n防御:{2}\n"
$ws.Range("I36").Value = "饱腹:{0}\n治疗:{1}\n防御:{2}\n"
$ws.Range("I37").Value = "饱腹:{0}\n治疗:{1}\n力量:{2}\n闪避:{3}\n"
$ws.Range("I38").Value = "饱腹:{0}\n治疗:{1}\n力量:{2}\n命中:{3}\n"
$ws.Range("I39").Value = "饱腹:{0}\n治疗:{1}\n防御:{2}\n速度:{3}\n"
$ws.Range("I40").Value = "饱腹:{0}\n幸运:{1}\n"
$ws.Range("I41").Value = "饱腹:{0}\n经验:{1}\n"
$ws.Range("I42").Value = "饱腹:{0}"

# Update the scroll position of the sheet view (topLeftCell A30 -> G30)
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 30
